# TRGP_sentiment.xlsx — append newer business-day rows to the sentiment table.
# Sheet1 previously ended at row 3662 (date serial 45835, 2025-06-27). This
# extends it through row 3716 (date serial 45915, 2025-09-15) with the
# Sentiment column (B) left at 0 for every new date, matching the existing
# pattern for the rest of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstNewRow = 3663
$lastExistingRow = $firstNewRow - 1

# Excel date-serial numbers for the new rows (business days, US holidays skipped)
$dateSerials = @(
    45838,45839,45840,45841,
    45845,45846,45847,45848,45849,
    45852,45853,45854,45855,45856,
    45859,45860,45861,45862,45863,
    45866,45867,45868,45869,45870,
    45873,45874,45875,45876,45877,
    45880,45881,45882,45883,45884,
    45887,45888,45889,45890,45891,
    45894,45895,45896,45897,45898,
    45902,45903,45904,45905,
    45908,45909,45910,45911,45912,
    45915
)

# Re-use the date number format already applied to column A so the new
# cells render the same way as the existing rows.
$dateFormat = $ws.Range("A$lastExistingRow").NumberFormat

for ($i = 0; $i -lt $dateSerials.Length; $i++) {
    $row = $firstNewRow + $i

    $aCell = $ws.Cells.Item($row, 1)
    $aCell.Value = $dateSerials[$i]
    $aCell.NumberFormat = $dateFormat

    $ws.Cells.Item($row, 2).Value = 0
}
